# ItemDetail TCER is updated
# Adds the "ItemDetail Service" TCER rows (R001-R004, two test blocks) and
# the "ItemDetail Dao" TCER rows to the "Item Details Service" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Block 1: ItemDetail Service (getAllItemDetails / findByItemId / ... ) ---

$ws.Range("A4").Value = 'R001'
$ws.Range("B4").Value = 'T001'
$ws.Range("C4").Value = 'getAllItemDetails()'
$ws.Range("D4").Value = 'List<ItemDetails> '
$ws.Range("F4").Value = 'collects all the items available in store'

$ws.Range("A5").Value = 'R001'
$ws.Range("B5").Value = 'T002'
$ws.Range("C5").Value = 'getAllItemDetails()'
$ws.Range("D5").Value = 'null'
$ws.Range("F5").Value = 'if there is no item present at store'

$ws.Range("A7").Value = 'R002'
$ws.Range("B7").Value = 'T001'
$ws.Range("C7").Value = 'findByItemId(String itemId)'
$ws.Range("D7").Value = 'ItemDetails item'
$ws.Range("F7").Value = 'Returns the ItemDetail object with the provided id'

$ws.Range("A8").Value = 'R002'
$ws.Range("B8").Value = 'T002'
$ws.Range("C8").Value = 'findByItemId(String itemId)'
$ws.Range("D8").Value = 'null'
$ws.Range("F8").Value = 'If item id deos not exist in store'

$ws.Range("A10").Value = 'R003'
$ws.Range("B10").Value = 'T001'
$ws.Range("C10").Value = 'findByItemId_AndAvailableQuantity(String Item Id , int availabkeQuantity)'
$ws.Range("D10").Value = 'ItemDetail item'
$ws.Range("F10").Value = 'returns Item detail if the provided item id consists of required quantity'

$ws.Range("A11").Value = 'R003'
$ws.Range("B11").Value = 'T002'
$ws.Range("C11").Value = 'findByItemId_AndAvailableQuantity(String Item Id , int availabkeQuantity)'
$ws.Range("D11").Value = 'null'
$ws.Range("F11").Value = 'if item with required quantity is not available'

$ws.Range("A13").Value = 'R004 '
$ws.Range("B13").Value = 'T001'
$ws.Range("C13").Value = 'updateRecord(String itemId,int quantity)'
$ws.Range("D13").Value = 'updated ItemDetail'
$ws.Range("F13").Value = 'returns the updated item id by reducing the mentioned quantity'

$ws.Range("A14").Value = 'R004 '
$ws.Range("B14").Value = 'T002'
$ws.Range("C14").Value = 'updateRecord(String itemId,int quantity)'
$ws.Range("D14").Value = 'null'
$ws.Range("F14").Value = 'if it cannot update the quantity mentioned'

# --- Section header: "ItemDetail Dao" (merged + formatted like the other section headers) ---

$ws.Range("A2:F2").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$ws.Range("A16:F16").Merge()
$ws.Range("A16").Value = 'ItemDetail Dao'

# --- Column header row for the ItemDetail Dao block (same style as the other header rows) ---

$ws.Range("A3:F3").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)

$ws.Range("A17").Value = 'Requirement #'
$ws.Range("B17").Value = 'TestCaseId'
$ws.Range("C17").Value = 'TestCase'
$ws.Range("D17").Value = 'Expected Result'
$ws.Range("E17").Value = 'Actual Result'
$ws.Range("F17").Value = 'Comment'

# --- Block 2: ItemDetail Dao (findAll / findById / findByItemIdAndAvailableQuantity / updateRecord) ---

$ws.Range("A18").Value = 'R001'
$ws.Range("B18").Value = 'T001'
$ws.Range("C18").Value = 'findAll()'
$ws.Range("D18").Value = 'List<ItemDetail>'
$ws.Range("F18").Value = 'returns the list of available items'

$ws.Range("A19").Value = 'R001'
$ws.Range("B19").Value = 'T002'
$ws.Range("C19").Value = 'findAll()'
$ws.Range("D19").Value = 'empty list'
$ws.Range("F19").Value = 'if there exist no element'

$ws.Range("A21").Value = 'R002'
$ws.Range("B21").Value = 'T001'
$ws.Range("C21").Value = 'findById(String itemId)'
$ws.Range("D21").Value = 'ItemDetail'
$ws.Range("F21").Value = 'if item with item id exists'

$ws.Range("A22").Value = 'R002'
$ws.Range("B22").Value = 'T002'
$ws.Range("C22").Value = 'findById(String itemId)'
$ws.Range("D22").Value = 'null'
$ws.Range("F22").Value = 'if item with item id does not exists'

$ws.Range("A24").Value = 'R003'
$ws.Range("B24").Value = 'T001'
$ws.Range("C24").Value = 'findByItemIdAndAvailableQuantity(String item id,int quantity)'
$ws.Range("D24").Value = 'ItemDetail'
$ws.Range("F24").Value = 'returns Item detail if the provided item id consists of required quantity'

$ws.Range("A25").Value = 'R003'
$ws.Range("B25").Value = 'T002'
$ws.Range("C25").Value = 'findByItemIdAndAvailableQuantity(String item id,int quantity)'
$ws.Range("D25").Value = 'null'
$ws.Range("F25").Value = 'if item with required quantity is not available'

$ws.Range("A27").Value = 'R004 '
$ws.Range("B27").Value = 'T001'
$ws.Range("C27").Value = 'updateRecord(String itemId,int quantity)'
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 'if item id updated by reducing the mentioned quantity'

$ws.Range("A28").Value = 'R004 '
$ws.Range("B28").Value = 'T002'
$ws.Range("C28").Value = 'updateRecord(String itemId,int quantity)'
$ws.Range("D28").Value = 0
$ws.Range("F28").Value = 'if it cannot update the quantity mentioned'

# --- Column widths for the now much wider TestCase / Comment columns ---

$ws.Columns.Item(3).ColumnWidth = 60.833333333333336
$ws.Columns.Item(6).ColumnWidth = 58.166666666666664

# --- View state: scroll/selection as left by the author ---

$ws.Select()
$ws.Range("F30").Select()
